$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.508.72"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.945.06"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "243.77"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("D7").Value = "58.39"
$ws.Range("E7").Value = "  -6.92%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.368"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").Value = "55.79"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "0.0838"
$ws.Range("E11").Value = "  +4.02%  "
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "0.829"
$ws.Range("E13").Value = "  -4.79%  "
$ws.Range("D14").Value = "21.53"
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("D15").Value = "2.227.61"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "13.62"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").Value = "5.27"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").Value = "1.931.01"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").Value = "36.417.86"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "69.85"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").Value = "0.0₃0872"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "229.86"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").Value = "5.01"
$ws.Range("E23").Value = "  -5.41%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "9.34"
$ws.Range("E27").Value = "  -5.14%  "
$ws.Range("D28").Value = "162.71"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").Value = "19.44"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  -9.29%  "
$ws.Range("D31").Value = "0.118"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.70"
$ws.Range("E33").Value = "  -4.86%  "
$ws.Range("D34").Value = "0.0633"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").Value = "6.24"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("D39").Value = "2.16"
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").Value = "0.0973"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "1.19"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.10"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("D47").Value = "1.352.71"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "87.92"
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("D49").Value = "7.22"
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("D50").Value = "2.82"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "45.76"
$ws.Range("E51").Value = "  +3.42%  "
